$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The refresh script re-pulled the group and a handful of fixtures
# that share an identical kickoff timestamp got re-ordered (their
# betting-odds/score/url payload in columns F:V swapped rows while the
# Indice/pais/torneio/temporada/data_partida columns A:E stayed put),
# plus three freshly played fixtures were appended at the bottom.
# ------------------------------------------------------------------

function Swap-RowPayload {
    param($rowA, $rowB)
    $va = $ws.Range("F$rowA`:V$rowA").Value()
    $vb = $ws.Range("F$rowB`:V$rowB").Value()
    $ws.Range("F$rowA`:V$rowA").Value = $vb
    $ws.Range("F$rowB`:V$rowB").Value = $va
}

# Rows 29/30/31 share the same timestamp and rotate: new29 <- old31,
# new30 <- old29, new31 <- old30.
$v29 = $ws.Range("F29:V29").Value()
$v30 = $ws.Range("F30:V30").Value()
$v31 = $ws.Range("F31:V31").Value()
$ws.Range("F29:V29").Value = $v31
$ws.Range("F30:V30").Value = $v29
$ws.Range("F31:V31").Value = $v30

# The remaining affected fixtures are simple pairwise swaps.
Swap-RowPayload 43 44
Swap-RowPayload 48 49
Swap-RowPayload 54 55
Swap-RowPayload 59 60
Swap-RowPayload 83 84
Swap-RowPayload 102 103

# ------------------------------------------------------------------
# Append the three newly scraped fixtures as rows 116-118, copying
# formatting from the last existing row (115) so the Indice (col A)
# and data_partida (col E) cells keep their original number formats.
# ------------------------------------------------------------------

function Add-MatchRow {
    param(
        $row, $indice, $pais, $torneio, $temporada, $dataPartida,
        $home, $homeGols, $away, $awayGols,
        $homeOpenOdds, $homeOpenDt, $homeCloseOdds, $homeCloseDt,
        $drawOpenOdds, $drawOpenDt, $drawCloseOdds, $drawCloseDt,
        $awayOpenOdds, $awayOpenDt, $awayCloseOdds, $awayCloseDt,
        $url
    )

    $ws.Range("A115:V115").Copy()
    $ws.Range("A$row`:V$row").PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $indice
    $ws.Cells.Item($row, 2).Value = $pais
    $ws.Cells.Item($row, 3).Value = $torneio
    $ws.Cells.Item($row, 4).Value = $temporada
    $ws.Cells.Item($row, 5).Value = $dataPartida
    $ws.Cells.Item($row, 6).Value = $home
    $ws.Cells.Item($row, 7).Value = $homeGols
    $ws.Cells.Item($row, 8).Value = $away
    $ws.Cells.Item($row, 9).Value = $awayGols
    $ws.Cells.Item($row, 10).Value = $homeOpenOdds
    $ws.Cells.Item($row, 11).Value = $homeOpenDt
    $ws.Cells.Item($row, 12).Value = $homeCloseOdds
    $ws.Cells.Item($row, 13).Value = $homeCloseDt
    $ws.Cells.Item($row, 14).Value = $drawOpenOdds
    $ws.Cells.Item($row, 15).Value = $drawOpenDt
    $ws.Cells.Item($row, 16).Value = $drawCloseOdds
    $ws.Cells.Item($row, 17).Value = $drawCloseDt
    $ws.Cells.Item($row, 18).Value = $awayOpenOdds
    $ws.Cells.Item($row, 19).Value = $awayOpenDt
    $ws.Cells.Item($row, 20).Value = $awayCloseOdds
    $ws.Cells.Item($row, 21).Value = $awayCloseDt
    $ws.Cells.Item($row, 22).Value = $url
}

Add-MatchRow 116 115 "czech-republic" "cfl-group-a" "2023-2024" 45241.42708333334 `
    "Motorlet Prague" 2 "FK Robstav" 3 `
    2.26 "10/11/2023 23:42" 2.3 "11/11/2023 10:07" `
    3.48 "10/11/2023 23:42" 3.56 "11/11/2023 10:07" `
    2.7 "10/11/2023 23:42" 2.71 "11/11/2023 10:07" `
    "https://www.betexplorer.com/football/czech-republic/cfl-group-a/motorlet-prague-fk-robstav/YZ32MBXA/"

Add-MatchRow 117 116 "czech-republic" "cfl-group-a" "2023-2024" 45241.58333333334 `
    "Domazlice" 5 "Admira Prague" 0 `
    1.31 "11/11/2023 10:43" 1.32 "11/11/2023 13:43" `
    5.19 "11/11/2023 10:43" 5.23 "11/11/2023 13:44" `
    6.95 "11/11/2023 10:43" 7.31 "11/11/2023 13:44" `
    "https://www.betexplorer.com/football/czech-republic/cfl-group-a/domazlice-admira-prague/0lKr2snd/"

Add-MatchRow 118 117 "czech-republic" "cfl-group-a" "2023-2024" 45241.58333333334 `
    "Taborsko akademie" 2 "Kraluv Dvur" 0 `
    1.48 "11/11/2023 10:43" 1.4 "11/11/2023 13:44" `
    4.47 "11/11/2023 10:43" 4.78 "11/11/2023 13:44" `
    4.94 "11/11/2023 10:43" 6.18 "11/11/2023 13:44" `
    "https://www.betexplorer.com/football/czech-republic/cfl-group-a/taborsko-akademie-kraluv-dvur/SfCUPpv3/"
